$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1: "Vaša VšZP" -> "TESTING"
$ws.Range("A1").Value = "TESTING"

# A4: value stays the same text "420987123452"
# (its shared-string index changes internally, content unchanged)
$ws.Range("A4").Value = "420987123452"

# Update the active selection from A5 to A2
$ws.Range("A2").Select()
